# "ATT: arrumado erros nas contas basicas"
#
# The "Caso de Teste" sheet (originally "Planilha1") lists manual test
# cases for calcular_juros_compostos(). This edit:
#   - renames the sheet
#   - replaces the three "insufficient values / empty tuple" rows with
#     three new, more specific "missing argument" test cases
#   - rewords the empty-tuple error message
#   - drops the trailing "value == 0" test rows (one of them, row 14,
#     is removed outright; the sheet now ends at row 13)
#   - keeps the remaining "string value" / "negative value" test rows,
#     which simply shift up as a result

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Caso de Teste"

$ws.Range("C4").Value = 'Não enviando o valor da capital'
$ws.Range("D4").Value = '20, 2'
$ws.Range("E4").Value = 'ValueError("Não é possível fazer a conta sem o valor da capital")'

$ws.Range("C5").Value = 'Não enviando o valor da taxa de juros'
$ws.Range("D5").Value = '1000, 2'
$ws.Range("E5").Value = 'ValueError("Não é possível fazer a conta sem o valor da taxa de juros")'

$ws.Range("C6").Value = 'Não enviando o valor do tempo'
$ws.Range("D6").Value = '1000, 30'
$ws.Range("E6").Value = 'ValueError("Não é possível fazer a conta sem o valor do tempo")'

$ws.Range("C7").Value = 'Testando envio de tupla vazia'
$ws.Range("D7").Value = '()'
$ws.Range("E7").Value = 'ValueError("Não é permitido uma tupla vazia")'

$ws.Range("C8").Value = 'Enviando uma string como valor da capital'
$ws.Range("D8").Value = '"oiii", 25, 2'
$ws.Range("E8").Value = 'ValueError("O valor do capital precisa ser int ou float")'

$ws.Range("C9").Value = 'Enviando uma string como valor dos juros'
$ws.Range("D9").Value = ' 2600, "oiii", 3'
$ws.Range("E9").Value = 'ValueError("O valor dos juros precisa ser int ou float")'

$ws.Range("C10").Value = 'Enviando uma string como valor do tempo'
$ws.Range("D10").Value = ' 2700, 4,"oiii"'
$ws.Range("E10").Value = 'ValueError("O valor do tempo precisa ser int ou float")'

$ws.Range("C11").Value = 'Enviando um valor da capital negativo'
$ws.Range("D11").Value = '-2000, 30, 3'
$ws.Range("E11").Value = 'ValueError("Não é permitido o valor da capital < que 0")'

$ws.Range("C12").Value = 'Enviando um valor dos juros negativo'
$ws.Range("D12").Value = '2000, -30, 3'
$ws.Range("E12").Value = 'ValueError("Não é permitido o valor dos juros < que 0")'

$ws.Range("C13").Value = 'Enviando um valor do tempo negativo'
$ws.Range("D13").Value = '2000, 30, -3'
$ws.Range("E13").Value = 'ValueError("Não é permitido o valor do tempo < que 0")'

# The "text" (leading-space) number format used to sit on D9 (the old
# "-2000, 30, 3" row); that row's text now lives at D11, so move the
# format along with it instead of leaving D9 looking like stale text.
$ws.Range("D10").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D11").NumberFormat = "@"

# The last row (the "tempo == 0" test) was dropped entirely.
$ws.Rows("14").Delete()

# Column E grew a little wider to fit the new longest error message.
$ws.Columns("E").ColumnWidth = 67.17

# Match the saved cell selection.
$ws.Range("F11").Select() | Out-Null
